$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (incl. number format) from A328 down through A329:A343 so new date cells match existing formatting
$ws.Range("A328").Copy($ws.Range("A329:A343"))

$ws.Range("A329").Value2 = 44403
$ws.Range("B329").Value2 = 0
$ws.Range("C329").Value2 = 0
$ws.Range("D329").Value2 = 0

$ws.Range("A330").Value2 = 44404
$ws.Range("B330").Value2 = 0
$ws.Range("C330").Value2 = 0
$ws.Range("D330").Value2 = 0

$ws.Range("A331").Value2 = 44405
$ws.Range("B331").Value2 = 0
$ws.Range("C331").Value2 = 0
$ws.Range("D331").Value2 = 0

$ws.Range("A332").Value2 = 44406
$ws.Range("B332").Value2 = 0
$ws.Range("C332").Value2 = 0
$ws.Range("D332").Value2 = 0

$ws.Range("A333").Value2 = 44407
$ws.Range("B333").Value2 = 0
$ws.Range("C333").Value2 = 0
$ws.Range("D333").Value2 = 0

$ws.Range("A334").Value2 = 44408
$ws.Range("B334").Value2 = 0
$ws.Range("C334").Value2 = 0
$ws.Range("D334").Value2 = 0

$ws.Range("A335").Value2 = 44409
$ws.Range("B335").Value2 = 1
$ws.Range("C335").Value2 = 1
$ws.Range("D335").Value2 = 40.79967360261118

$ws.Range("A336").Value2 = 44410
$ws.Range("B336").Value2 = 0
$ws.Range("C336").Value2 = 1
$ws.Range("D336").Value2 = 40.79967360261118

$ws.Range("A337").Value2 = 44411
$ws.Range("B337").Value2 = 1
$ws.Range("C337").Value2 = 2
$ws.Range("D337").Value2 = 81.59934720522236

$ws.Range("A338").Value2 = 44412
$ws.Range("B338").Value2 = 0
$ws.Range("C338").Value2 = 2
$ws.Range("D338").Value2 = 81.59934720522236

$ws.Range("A339").Value2 = 44413
$ws.Range("B339").Value2 = 1
$ws.Range("C339").Value2 = 3
$ws.Range("D339").Value2 = 122.3990208078335

$ws.Range("A340").Value2 = 44414
$ws.Range("B340").Value2 = 1
$ws.Range("C340").Value2 = 4
$ws.Range("D340").Value2 = 163.1986944104447

$ws.Range("A341").Value2 = 44415
$ws.Range("B341").Value2 = 6
$ws.Range("C341").Value2 = 10
$ws.Range("D341").Value2 = 407.9967360261118

$ws.Range("A342").Value2 = 44416
$ws.Range("B342").Value2 = 3
$ws.Range("C342").Value2 = 12
$ws.Range("D342").Value2 = 489.5960832313341

$ws.Range("A343").Value2 = 44417
$ws.Range("B343").Value2 = 1
$ws.Range("C343").Value2 = 13
$ws.Range("D343").Value2 = 530.3957568339454